$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Ağırlık" value for row 19 (F3)
$ws.Range("D19").Value = "AĞIR"

# Add a new row of data: G1, 10, İÇECEK, AĞIR
$ws.Range("A20").Value = "G1"
$ws.Range("B20").Value = 10
$ws.Range("C20").Value = "İÇECEK"
$ws.Range("D20").Value = "AĞIR"

# Update the current selection to match the new active cell
$ws.Range("E20").Select()
